$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.131980458546578
$ws.Range("F3").Value = 2.146819819230537
$ws.Range("F4").Value = 2.526153402358124
$ws.Range("F5").Value = 19
$ws.Range("F7").Value = 2.179755441148193
$ws.Range("F8").Value = 0.78
$ws.Range("F9").Value = 1.296185451270793
